# Adjust New Dummy Database
$wb = $excel.ActiveWorkbook

$wsSiteList   = $wb.Worksheets.Item("Site List")
$wsInsertRing = $wb.Worksheets.Item("Insert Ring")

# --- Update "Insert Ring" data (row 18 / row 19) ---
$wsInsertRing.Range("L18").Value = "11PDG0282"

$wsInsertRing.Range("K19").Value = "11PDG0262"
$wsInsertRing.Range("M19").Value = "11PDG0262-11PDG0282"
$wsInsertRing.Range("R19").Value = "11PDG0282-02RKB007"
$wsInsertRing.Range("S19").Value = 2100
$wsInsertRing.Range("T19").Value = 1000
$wsInsertRing.Range("U19").Value = 3100
$wsInsertRing.Range("V19").Value = "Segment Insert"

# --- Switch the active/selected sheet & selection ---
# Before: "Site List" tab selected, selection G9; "Insert Ring" selection N11
# After:  "Insert Ring" tab selected, selection D31; "Site List" selection G14
$wsSiteList.Activate()
$wsSiteList.Range("G14").Select()

$wsInsertRing.Activate()
$wsInsertRing.Range("D31").Select()

Write-Output "done"
